# Add a new "dob" entity row to the DB schema table (Table3), between the
# "aadharNo" row and the "streetAdress" row, and give the "Course Name (PG)"
# row a distinct (bottom-less) border once the new row pushes it to row 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new sheet row at row 8 (the table's 7th data row), shifting the
# rows below it (and their formatting) down by one.
$ws.Rows("8").Insert()

# Grow the table/autofilter range to match the newly-inserted row.
$lo.Resize($ws.Range("A1:D23"))

# Populate the new row: Entity Name / Data Type / Label Name (Remarks stays blank).
$ws.Range("A8").Value = "dob"
$ws.Range("B8").Value = "date"
$ws.Range("C8").Value = "Date of birth"

# Give the new row the same (bordered) formatting as the rest of the table
# by copying formats from the row directly above it.
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)

# The "Course Name (PG)" row has shifted down to row 22 - mark it with a
# thin border on the left/top/right but no bottom border.
$rng22 = $ws.Range("A22:D22")
$rng22.Borders.LineStyle = 1
$rng22.Borders.Item(9).LineStyle = 0

# Match the saved selection from the source workbook.
$ws.Range("G15").Select()
